# lesson 4 examples, exercies and slide
#
# 1) Notes Master "Date Placeholder" auto date field was re-cached by
#    PowerPoint from 2/29/2024 to 3/11/2024 (the field keeps updating to
#    "today" whenever the deck is opened/saved).
# 2) Slide 1's presenter-name placeholder changed from the Hungarian
#    author "Jozsef Gal" (hu-HU) to "Janos Stefan" (en-US).

$p = $ppt.ActivePresentation

# --- 1. Notes Master date/time footer field -------------------------------
$nm = $p.NotesMaster
$nm.HeadersFooters.DateAndTime.Text = "3/11/2024"

# --- 2. Slide 1 presenter name placeholder ---------------------------------
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)
$sh.TextFrame.TextRange.Text = "Janos Stefan"
$sh.TextFrame2.TextRange.LanguageID = "en-US"
